# Auto-applied cell updates per commit diff (Typhon_Profits.xlsx scheduled market-data refresh)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 73277.86
$ws.Range("I137").Value = 1748.5
$ws.Range("K137").Value = 5245.5
$ws.Range("M137").Value = -2695.5
$ws.Range("H138").Value = 1833.629
$ws.Range("I138").Value = 1206.2258
$ws.Range("J138").Value = 2461.0322
$ws.Range("K138").Value = 3618.6774
$ws.Range("L138").Value = 7383.096600000001
$ws.Range("M138").Value = 1521.3226
$ws.Range("N138").Value = -17663.0966
$ws.Range("H141").Value = 2811.111
$ws.Range("I141").Value = 2155.7144
$ws.Range("J141").Value = 5105
$ws.Range("K141").Value = 6467.1432
$ws.Range("L141").Value = 15315
$ws.Range("M141").Value = -1287.1432
$ws.Range("N141").Value = -25675
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19536.316
$ws.Range("I32").Value = 19876.25
$ws.Range("K32").Value = 19876.25
$ws.Range("M32").Value = -19589.25
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H61").Value = 2327.9
$ws.Range("I61").Value = 1929.6
$ws.Range("K61").Value = 1929.6
$ws.Range("M61").Value = -1717.6
$ws.Range("H74").Value = 52634504
$ws.Range("I74").Value = 66670036
$ws.Range("J74").Value = 1253.5
$ws.Range("K74").Value = 66670036
$ws.Range("L74").Value = 1253.5
$ws.Range("M74").Value = -66669162
$ws.Range("N74").Value = -3001.5
$ws.Range("H77").Value = 52634504
$ws.Range("I77").Value = 66670036
$ws.Range("J77").Value = 1253.5
$ws.Range("K77").Value = 333350180
$ws.Range("L77").Value = 6267.5
$ws.Range("M77").Value = -333345812
$ws.Range("N77").Value = -15003.5
$ws.Range("H110").Value = 690
$ws.Range("I110").Value = 690
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 690
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 1355
$ws.Range("N110").ClearContents()
$ws.Range("H132").Value = 15731.027
$ws.Range("I132").Value = 1645.4584
$ws.Range("K132").Value = 4936.3752
$ws.Range("M132").Value = -2406.3752
$ws.Range("H136").Value = 2327.9
$ws.Range("I136").Value = 1929.6
$ws.Range("K136").Value = 5788.799999999999
$ws.Range("M136").Value = -3238.799999999999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1719.091
$ws.Range("I99").Value = 985
$ws.Range("K99").Value = 985
$ws.Range("M99").Value = 513
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 7813400
$ws.Range("I105").Value = 13889516
$ws.Range("K105").Value = 13889516
$ws.Range("M105").Value = -13887769
$ws.Range("H132").Value = 16918.885
$ws.Range("I132").Value = 19074.104
$ws.Range("J132").Value = 6502
$ws.Range("K132").Value = 57222.312
$ws.Range("L132").Value = 19506
$ws.Range("M132").Value = -54692.312
$ws.Range("N132").Value = -24566
$ws.Range("H134").Value = 1017.46344
$ws.Range("I134").Value = 918.7727
$ws.Range("K134").Value = 2756.3181
$ws.Range("M134").Value = -221.3181
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1329.5714
$ws.Range("I5").Value = 717.8333
$ws.Range("K5").Value = 2153.4999
$ws.Range("M5").Value = -2041.4999
$ws.Range("H18").Value = 206.66667
$ws.Range("I18").Value = 158
$ws.Range("J18").Value = 450
$ws.Range("K18").Value = 474
$ws.Range("L18").Value = 1350
$ws.Range("M18").Value = -305
$ws.Range("N18").Value = -1688
$ws.Range("H36").Value = 101912.22
$ws.Range("J36").Value = 130344
$ws.Range("L36").Value = 391032
$ws.Range("N36").Value = -391370
$ws.Range("H113").Value = 4261.185
$ws.Range("I113").Value = 7533.857
$ws.Range("J113").Value = 736.7692
$ws.Range("K113").Value = 22601.571
$ws.Range("L113").Value = 2210.3076
$ws.Range("M113").Value = -20431.571
$ws.Range("N113").Value = -6550.3076
$ws.Range("H122").Value = 661.625
$ws.Range("I122").Value = 275
$ws.Range("J122").Value = 1048.25
$ws.Range("K122").Value = 2475
$ws.Range("L122").Value = 9434.25
$ws.Range("M122").Value = -25
$ws.Range("N122").Value = -14334.25
$ws.Range("H131").Value = 725.5
$ws.Range("J131").Value = 725.5
$ws.Range("L131").Value = 2176.5
$ws.Range("N131").Value = -12256.5
$ws.Range("H135").Value = 1329.5714
$ws.Range("I135").Value = 717.8333
$ws.Range("K135").Value = 6460.4997
$ws.Range("M135").Value = -3925.4997
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2140.7036
$ws.Range("J113").Value = 2707.25
$ws.Range("L113").Value = 2707.25
$ws.Range("N113").Value = -7047.25
$ws.Range("H132").Value = 46360.656
$ws.Range("I132").Value = 47582.61
$ws.Range("J132").Value = 44018.582
$ws.Range("K132").Value = 142747.83
$ws.Range("L132").Value = 132055.746
$ws.Range("M132").Value = -140217.83
$ws.Range("N132").Value = -137115.746
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()
$ws.Range("H61").Value = 2788.1516
$ws.Range("J61").Value = 5991
$ws.Range("L61").Value = 5991
$ws.Range("N61").Value = -6395
$ws.Range("H68").Value = 2470.2666
$ws.Range("I68").Value = 2205.9
$ws.Range("J68").Value = 2999
$ws.Range("K68").Value = 2205.9
$ws.Range("L68").Value = 2999
$ws.Range("M68").Value = -1456.9
$ws.Range("N68").Value = -4497
$ws.Range("H71").Value = 2470.2666
$ws.Range("I71").Value = 2205.9
$ws.Range("J71").Value = 2999
$ws.Range("K71").Value = 11029.5
$ws.Range("L71").Value = 14995
$ws.Range("M71").Value = -7285.5
$ws.Range("N71").Value = -22483
$ws.Range("H113").Value = 2788.1516
$ws.Range("J113").Value = 5991
$ws.Range("L113").Value = 5991
$ws.Range("N113").Value = -10331
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1196.1578
$ws.Range("I132").Value = 770.375
$ws.Range("J132").Value = 3467
$ws.Range("K132").Value = 2311.125
$ws.Range("L132").Value = 10401
$ws.Range("M132").Value = 218.875
$ws.Range("N132").Value = -15461
